$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 119.681816
$ws.Range("I11").Value = 119.681816
$ws.Range("K11").Value = 119.681816
$ws.Range("M11").Value = 20.318184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 10001512
$ws.Range("I33").Value = 15625272
$ws.Range("J33").Value = 3715.6667
$ws.Range("K33").Value = 15625272
$ws.Range("L33").Value = 3715.6667
$ws.Range("M33").Value = -15625043
$ws.Range("N33").Value = -4173.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 80799.89
$ws.Range("J75").Value = 80799.89
$ws.Range("L75").Value = 80799.89
$ws.Range("N75").Value = -82671.89

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 80799.89
$ws.Range("J78").Value = 80799.89
$ws.Range("L78").Value = 242399.67
$ws.Range("N78").Value = -251759.67

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1786.3334
$ws.Range("J88").Value = 1755.6
$ws.Range("L88").Value = 1755.6
$ws.Range("N88").Value = -2567.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1786.3334
$ws.Range("J91").Value = 1755.6
$ws.Range("L91").Value = 1755.6
$ws.Range("N91").Value = -4563.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3216.2415
$ws.Range("I98").Value = 3270.92
$ws.Range("K98").Value = 3270.92
$ws.Range("M98").Value = -1772.92

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3216.2415
$ws.Range("I122").Value = 3270.92
$ws.Range("K122").Value = 9812.76
$ws.Range("M122").Value = -7362.76

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 456207.47
$ws.Range("I127").Value = 456207.47
$ws.Range("K127").Value = 1368622.41
$ws.Range("M127").Value = -1363662.41

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3645.0557
$ws.Range("I141").Value = 3766.6155
$ws.Range("K141").Value = 11299.8465
$ws.Range("M141").Value = -6119.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 291.58334
$ws.Range("I97").Value = 291.58334
$ws.Range("K97").Value = 291.58334
$ws.Range("M97").Value = 204.41666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 26749.5
$ws.Range("J112").Value = 26749.5
$ws.Range("L112").Value = 26749.5
$ws.Range("N112").Value = -29703.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 35000
$ws.Range("J114").Value = 35000
$ws.Range("L114").Value = 35000
$ws.Range("N114").Value = -43678

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4732614.5
$ws.Range("I132").Value = 2175.4092
$ws.Range("J132").Value = 22077556
$ws.Range("K132").Value = 6526.2276
$ws.Range("L132").Value = 66232668
$ws.Range("M132").Value = -3996.2276
$ws.Range("N132").Value = -66237728

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3927.25
$ws.Range("I99").Value = 8010
$ws.Range("J99").Value = 2566.3333
$ws.Range("K99").Value = 8010
$ws.Range("L99").Value = 2566.3333
$ws.Range("M99").Value = -6512
$ws.Range("N99").Value = -5562.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9662.621999999999
$ws.Range("I58").Value = 3446.9211
$ws.Range("K58").Value = 3446.9211
$ws.Range("M58").Value = -3243.9211

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11658.077
$ws.Range("I86").Value = 12416
$ws.Range("J86").Value = 9131.666999999999
$ws.Range("K86").Value = 12416
$ws.Range("L86").Value = 9131.666999999999
$ws.Range("M86").Value = -11293
$ws.Range("N86").Value = -11377.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 11658.077
$ws.Range("I89").Value = 12416
$ws.Range("J89").Value = 9131.666999999999
$ws.Range("K89").Value = 62080
$ws.Range("L89").Value = 45658.335
$ws.Range("M89").Value = -56464
$ws.Range("N89").Value = -56890.335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6158
$ws.Range("I99").Value = 2562
$ws.Range("J99").Value = 14249
$ws.Range("K99").Value = 2562
$ws.Range("L99").Value = 14249
$ws.Range("M99").Value = -1064
$ws.Range("N99").Value = -17245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6158
$ws.Range("I126").Value = 2562
$ws.Range("J126").Value = 14249
$ws.Range("K126").Value = 7686
$ws.Range("L126").Value = 42747
$ws.Range("M126").Value = -5216
$ws.Range("N126").Value = -47687

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 9662.621999999999
$ws.Range("I136").Value = 3446.9211
$ws.Range("K136").Value = 10340.7633
$ws.Range("M136").Value = -7790.763300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8298591.5
$ws.Range("I4").Value = 13418333
$ws.Range("J4").Value = 1716066.9
$ws.Range("K4").Value = 40254999
$ws.Range("L4").Value = 5148200.699999999
$ws.Range("M4").Value = -40254887
$ws.Range("N4").Value = -5148424.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4583.375
$ws.Range("J39").Value = 5850
$ws.Range("L39").Value = 17550
$ws.Range("N39").Value = -18138

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1474.0303
$ws.Range("I131").Value = 1038.3334
$ws.Range("J131").Value = 1502.1398
$ws.Range("K131").Value = 3115.0002
$ws.Range("L131").Value = 4506.4194
$ws.Range("M131").Value = 1924.9998
$ws.Range("N131").Value = -14586.4194

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3457.875
$ws.Range("I80").Value = 3544
$ws.Range("J80").Value = 3314.3333
$ws.Range("K80").Value = 3544
$ws.Range("L80").Value = 3314.3333
$ws.Range("M80").Value = -2546
$ws.Range("N80").Value = -5310.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3457.875
$ws.Range("I83").Value = 3544
$ws.Range("J83").Value = 3314.3333
$ws.Range("K83").Value = 17720
$ws.Range("L83").Value = 16571.6665
$ws.Range("M83").Value = -12728
$ws.Range("N83").Value = -26555.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2102.6667
$ws.Range("I97").Value = 2500
$ws.Range("K97").Value = 2500
$ws.Range("M97").Value = -2004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -46134

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 84210
$ws.Range("J54").Value = 84210
$ws.Range("L54").Value = 84210
$ws.Range("N54").Value = -85498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2515.875
$ws.Range("I82").Value = 3134.9285
$ws.Range("K82").Value = 3134.9285
$ws.Range("M82").Value = -2773.9285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2515.875
$ws.Range("I85").Value = 3134.9285
$ws.Range("K85").Value = 3134.9285
$ws.Range("M85").Value = -1886.9285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 16666.666
$ws.Range("J96").Value = 16666.666
$ws.Range("L96").Value = 16666.666
$ws.Range("N96").Value = -22158.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3017.4736
$ws.Range("I100").Value = 2636.1667
$ws.Range("J100").Value = 3671.1428
$ws.Range("K100").Value = 2636.1667
$ws.Range("L100").Value = 3671.1428
$ws.Range("M100").Value = -2095.1667
$ws.Range("N100").Value = -4753.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10512.75
$ws.Range("J122").Value = 4850.5
$ws.Range("L122").Value = 14551.5
$ws.Range("N122").Value = -19451.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1749.8334
$ws.Range("J96").Value = 1799.8
$ws.Range("L96").Value = 1799.8
$ws.Range("N96").Value = -4545.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 575
$ws.Range("I100").Value = 466.66666
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 933.33332
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -392.33332
$ws.Range("N100").Value = -2882

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 30615
$ws.Range("J105").Value = 30615
$ws.Range("L105").Value = 30615
$ws.Range("N105").Value = -37603

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2164.5386
$ws.Range("I122").Value = 1930
$ws.Range("J122").Value = 2692.25
$ws.Range("K122").Value = 5790
$ws.Range("L122").Value = 8076.75
$ws.Range("M122").Value = -3340
$ws.Range("N122").Value = -12976.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 91376.92
$ws.Range("J140").Value = 93158.336
$ws.Range("L140").Value = 93158.336
$ws.Range("N140").Value = -103518.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 96250
$ws.Range("I141").Value = 95000
$ws.Range("J141").Value = 96666.664
$ws.Range("K141").Value = 95000
$ws.Range("L141").Value = 96666.664
$ws.Range("M141").Value = -89820
$ws.Range("N141").Value = -107026.664
